# "Se actualiza los datos de pagos"
#
# The three data-driven sheets (CronogramaPagosVehicular, HistorialPagos,
# DescargaHistorialPagos) each hold one row of login/vehicle test data in a
# table (A2:D2 = numeroUsuario / contrasena / placa / vigencia). This commit
# refreshes that row to a single new data set on every sheet:
#   numeroUsuario = 72934725
#   contrasena    = Rimac2020 (unchanged)
#   placa         = XFN-363
#   vigencia      = 06/11/2025
#
# It also moves the active/selected tab from the 1st sheet to the 3rd sheet
# (DescargaHistorialPagos), and updates the selected cell on sheets 2 and 3
# to D2.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: CronogramaPagosVehicular -------------------------------------
$ws1 = $wb.Worksheets.Item(1)
# Leading apostrophes keep these as text (matching the original t="s" cells)
# instead of letting Excel auto-coerce them to a number / date.
$ws1.Range("A2").Value = "'72934725"
$ws1.Range("C2").Value = "'XFN-363"
$ws1.Range("D2").Value = "'06/11/2025"

# --- Sheet 2: HistorialPagos -------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "'72934725"
$ws2.Range("C2").Value = "'XFN-363"
$ws2.Range("D2").Value = "'06/11/2025"
$ws2.Range("D2").Select() | Out-Null

# --- Sheet 3: DescargaHistorialPagos ----------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "'72934725"
$ws3.Range("C2").Value = "'XFN-363"
$ws3.Range("D2").Value = "'06/11/2025"
$ws3.Columns.Item(4).ColumnWidth = 11

# Make DescargaHistorialPagos the active tab (it becomes the last-selected
# sheet, so select it - and its target cell - last).
$ws3.Range("D2").Select() | Out-Null
